{"js": "// Replace each two-digit-multiplication expression in the table with its\n// updated counterpart. Old values are unique across the document, so a\n// direct search-and-replace per pair is safe and order-independent.\nconst replacements = [\n  [\"43\u00d763=2709\", \"70\u00d755=3850\"],\n  [\"14\u00d732=448\", \"87\u00d740=3480\"],\n  [\"54\u00d721=1134\", \"16\u00d714=224\"],\n  [\"50\u00d748=2400\", \"96\u00d722=2112\"],\n  [\"62\u00d769=4278\", \"64\u00d759=3776\"],\n  [\"12\u00d782=984\", \"21\u00d740=840\"],\n  [\"73\u00d756=4088\", \"74\u00d752=3848\"],\n  [\"46\u00d753=2438\", \"41\u00d787=3567\"],\n  [\"50\u00d793=4650\", \"15\u00d769=1035\"],\n  [\"34\u00d716=544\", \"80\u00d741=3280\"],\n  [\"42\u00d751=2142\", \"19\u00d791=1729\"],\n  [\"12\u00d717=204\", \"53\u00d712=636\"],\n  [\"80\u00d753=4240\", \"53\u00d766=3498\"],\n  [\"78\u00d769=5382\", \"70\u00d755=3850\"],\n  [\"45\u00d767=3015\", \"90\u00d752=4680\"],\n  [\"20\u00d735=700\", \"96\u00d784=8064\"],\n  [\"29\u00d748=1392\", \"36\u00d773=2628\"],\n  [\"41\u00d746=1886\", \"40\u00d787=3480\"],\n  [\"75\u00d713=975\", \"18\u00d733=594\"],\n  [\"97\u00d756=5432\", \"13\u00d775=975\"],\n  [\"69\u00d782=5658\", \"93\u00d759=5487\"],\n  [\"33\u00d729=957\", \"90\u00d759=5310\"],\n  [\"40\u00d772=2880\", \"31\u00d787=2697\"],\n  [\"95\u00d730=2850\", \"81\u00d713=1053\"],\n  [\"58\u00d730=1740\", \"61\u00d776=4636\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-multiplication expression in the table with its\n# updated counterpart. Old values are unique across the document, so a\n# direct Find/Replace per pair (wdReplaceAll) is safe and order-independent.\n$pairs = @(\n    @(\"43\u00d763=2709\", \"70\u00d755=3850\"),\n    @(\"14\u00d732=448\", \"87\u00d740=3480\"),\n    @(\"54\u00d721=1134\", \"16\u00d714=224\"),\n    @(\"50\u00d748=2400\", \"96\u00d722=2112\"),\n    @(\"62\u00d769=4278\", \"64\u00d759=3776\"),\n    @(\"12\u00d782=984\", \"21\u00d740=840\"),\n    @(\"73\u00d756=4088\", \"74\u00d752=3848\"),\n    @(\"46\u00d753=2438\", \"41\u00d787=3567\"),\n    @(\"50\u00d793=4650\", \"15\u00d769=1035\"),\n    @(\"34\u00d716=544\", \"80\u00d741=3280\"),\n    @(\"42\u00d751=2142\", \"19\u00d791=1729\"),\n    @(\"12\u00d717=204\", \"53\u00d712=636\"),\n    @(\"80\u00d753=4240\", \"53\u00d766=3498\"),\n    @(\"78\u00d769=5382\", \"70\u00d755=3850\"),\n    @(\"45\u00d767=3015\", \"90\u00d752=4680\"),\n    @(\"20\u00d735=700\", \"96\u00d784=8064\"),\n    @(\"29\u00d748=1392\", \"36\u00d773=2628\"),\n    @(\"41\u00d746=1886\", \"40\u00d787=3480\"),\n    @(\"75\u00d713=975\", \"18\u00d733=594\"),\n    @(\"97\u00d756=5432\", \"13\u00d775=975\"),\n    @(\"69\u00d782=5658\", \"93\u00d759=5487\"),\n    @(\"33\u00d729=957\", \"90\u00d759=5310\"),\n    @(\"40\u00d772=2880\", \"31\u00d787=2697\"),\n    @(\"95\u00d730=2850\", \"81\u00d713=1053\"),\n    @(\"58\u00d730=1740\", \"61\u00d776=4636\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nWrite-Output \"done\"\n"}
